$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.0782454872199231
$ws.Range("B3").Value = 23.13929419043147
$ws.Range("B4").Value = 151.7000066755207
$ws.Range("B5").Value = 0.1862698743118686
$ws.Range("B6").Value = 2.550013287552345
$ws.Range("B7").Value = 0.7589387427737242
$ws.Range("B9").Value = 1.830845269518987
$ws.Range("B11").Value = 914.0617569498547
$ws.Range("B12").Value = 0.9547297179517178
$ws.Range("B13").Value = 0.9393391234223979
$ws.Range("B14").Value = 2.438758919211486
$ws.Range("B15").Value = 0.914490606252917
$ws.Range("B16").Value = 0.07821515985127751
$ws.Range("B19").Value = 0.4594743153861449
$ws.Range("B20").Value = 0.1275017474602679
$ws.Range("B21").Value = 0.1265271498996089
$ws.Range("B22").Value = 5.107268771530819
$ws.Range("B23").Value = 0.2382463473900343
$ws.Range("B24").Value = 0.4067438140266562
$ws.Range("B25").Value = 28.30467215265686
$ws.Range("B26").Value = 29.29985665220562
$ws.Range("B27").Value = 0.08220955371778636
$ws.Range("B28").Value = 0.4548255380773952
$ws.Range("B29").Value = 1.708027906264028
$ws.Range("B30").Value = 3.595510537878962
$ws.Range("B31").Value = 12.06648340382061
$ws.Range("B32").Value = 33.16550990136131
$ws.Range("B33").Value = 73852.54308902436
$ws.Range("B34").Value = 9.001542361301979
$ws.Range("B35").Value = 87.63091361475982
$ws.Range("B36").Value = 152.2856605377513
